# Adding a SOAP WebService sheet ("soapMimeHeaders") to the Spring Boot / SAAJ
# test-cases workbook: a new worksheet listing the SOAP MIME headers
# (Content-Type / application/xml and SOAPAction / <url>), with the SOAPAction
# value also turned into a clickable hyperlink.

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet so it lands at the end
# of the tab strip (and becomes the active sheet, like a freshly inserted tab
# in Excel).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "soapMimeHeaders"

# Header/value pairs mirroring the other "Header"/"Parameters" sheets.
$ws.Range("A1").Value = "Content-Type"
$ws.Range("B1").Value = "application/xml"
$ws.Range("A2").Value = "SOAPAction"
$ws.Range("B2").Value = "http://application.com/soap/products/addProductRequest"

# Turn the SOAPAction value into a real hyperlink (adds the Hyperlink style +
# font, plus the external relationship).
$ws.Hyperlinks.Add($ws.Range("B2"), "http://application.com/soap/products/addProductRequest") | Out-Null

# Leave the selection on the hyperlink cell, matching the saved selection in
# the new sheet.
$ws.Range("B2").Select() | Out-Null
